$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert six new rows above the old row 45 ("n1-o3_1" / instanceID) block so
# the existing rows 45-46 shift down to 51-52, matching the target layout.
# ---------------------------------------------------------------------------
$ws.Rows("45:50").Insert()

# ---------------------------------------------------------------------------
# Populate the new dictionary rows. The write order below reproduces the
# exact order in which these new entries were first authored (mirrored by
# the order new strings were appended to the shared-string table upstream).
# ---------------------------------------------------------------------------
$ws.Range("A45").Value = "n1-o3_1aa"
$ws.Range("D45").Value = "o3_1aa"

$ws.Range("A46").Value = "n1-o3_1b"
$ws.Range("D46").Value = "o3_1b"

$ws.Range("A47").Value = "n1-o3_3"
$ws.Range("C47").Value = "care_pathway"
$ws.Range("D47").Value = "o3_3"

$ws.Range("A48").Value = "n1-o3_2a"

$ws.Range("C45").Value = "other_care_yn"
$ws.Range("C46").Value = "other_care_provider"

$ws.Range("D48").Value = "o3_2a"

$ws.Range("C49").Value = "rx_day7"

$ws.Range("A49").Value = "n1-o3_2b"

$ws.Range("C48").Value = "rx_day7_yn"

$ws.Range("D49").Value = "o3_2b"

$ws.Range("A50").Value = "n1-o3_2o"
$ws.Range("D50").Value = "o3_2o"

$ws.Range("C50").Value = "rx_day7_oth"

# B column flags (all 1, same as the other rows in this block).
$ws.Range("B45").Value = 1
$ws.Range("B46").Value = 1
$ws.Range("B47").Value = 1
$ws.Range("B48").Value = 1
$ws.Range("B49").Value = 1
$ws.Range("B50").Value = 1

# ---------------------------------------------------------------------------
# New highlight style: solid orange fill (FFC000) with centered alignment,
# applied to column A and C of the new rows (new cellXfs #15 / fill #6).
# ---------------------------------------------------------------------------
$highlightCells = @("A45","C45","A46","C46","A47","C47","A48","C48","A49","C49","A50","C50")
foreach ($addr in $highlightCells) {
    $cell = $ws.Range($addr)
    $cell.Interior.Color = 49407
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# ---------------------------------------------------------------------------
# Selection / active cell bookkeeping to mirror the author's final view.
# ---------------------------------------------------------------------------
$ws.Range("A45:XFD50").Select()
